$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.633.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.968.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.43%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.63"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.89%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.375"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.11%  "

$ws.Range("E12").Value = "  -0.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.844"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.256.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.971.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.492.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0886"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.73%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "166.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.123"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0645"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("E38").Value = "  -1.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0979"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.61%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.43%  "

$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0212"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.355.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.17%  "
